$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("B2").Value = 7191
$ws.Range("C3").Value = 170731
$ws.Range("C4").Value = 161555
$ws.Range("C8").Value = 65.75
